$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The refreshed crypto-price/volume figures below are stored as plain text
# (e.g. "303.75", "5.65%") just like the original cells, so for each one we
# temporarily switch the cell to Text format before writing the value - this
# stops Excel from auto-parsing the string into a number/percentage - and
# then restore the original "General" format.
$updates = @(
    @{ Cell = "D2"; Value = "303.75" },
    @{ Cell = "E2"; Value = "5.65%" },
    @{ Cell = "D3"; Value = "32.21" },
    @{ Cell = "E3"; Value = "10.54%" },
    @{ Cell = "D4"; Value = "5.272" },
    @{ Cell = "E4"; Value = "1.44%" },
    @{ Cell = "D5"; Value = "0.07485" },
    @{ Cell = "E5"; Value = "7.46%" },
    @{ Cell = "D6"; Value = "7.856" },
    @{ Cell = "E6"; Value = "5.69%" },
    @{ Cell = "E7"; Value = "7.16%" },
    @{ Cell = "D8"; Value = "1.505" },
    @{ Cell = "E8"; Value = "7.26%" },
    @{ Cell = "D9"; Value = "0.9203" },
    @{ Cell = "E9"; Value = "2.01%" },
    @{ Cell = "D10"; Value = "0.1690" },
    @{ Cell = "E10"; Value = "5.05%" },
    @{ Cell = "D11"; Value = "0.07996" },
    @{ Cell = "E11"; Value = "4.80%" },
    @{ Cell = "D12"; Value = "0.08054" },
    @{ Cell = "E12"; Value = "3.68%" },
    @{ Cell = "E13"; Value = "2.63%" },
    @{ Cell = "D14"; Value = "0.09901" },
    @{ Cell = "E14"; Value = "9.92%" },
    @{ Cell = "D15"; Value = "0.001504" },
    @{ Cell = "E15"; Value = "-6.01%" },
    @{ Cell = "D16"; Value = "0.04608" },
    @{ Cell = "E16"; Value = "1.77%" },
    @{ Cell = "D17"; Value = "0.006566" },
    @{ Cell = "E17"; Value = "0.45%" },
    @{ Cell = "E18"; Value = "-0.11%" },
    @{ Cell = "D19"; Value = "2.230" },
    @{ Cell = "E19"; Value = "0.01%" },
    @{ Cell = "D21"; Value = "0.1343" },
    @{ Cell = "E21"; Value = "0.41%" },
    @{ Cell = "D22"; Value = "4.486" },
    @{ Cell = "E22"; Value = "10.79%" },
    @{ Cell = "E23"; Value = "1.36%" },
    @{ Cell = "E24"; Value = "0.51%" },
    @{ Cell = "D25"; Value = "0.004449" },
    @{ Cell = "E25"; Value = "7.36%" },
    @{ Cell = "D26"; Value = "0.0001398" },
    @{ Cell = "E26"; Value = "19.53%" },
    @{ Cell = "E27"; Value = "6.51%" },
    @{ Cell = "D39"; Value = "0.01721" },
    @{ Cell = "E39"; Value = "2,556.56%" },
    @{ Cell = "D40"; Value = "0.04497" },
    @{ Cell = "E40"; Value = "2.90%" },
    @{ Cell = "D41"; Value = "0.007155" },
    @{ Cell = "E41"; Value = "3.02%" },
    @{ Cell = "D42"; Value = "0.1349" },
    @{ Cell = "E42"; Value = "8.43%" },
    @{ Cell = "D43"; Value = "0.002236" },
    @{ Cell = "E43"; Value = "8.09%" },
    @{ Cell = "D44"; Value = "0.01279" },
    @{ Cell = "E44"; Value = "10.01%" },
    @{ Cell = "D45"; Value = "0.00006164" },
    @{ Cell = "E45"; Value = "5.79%" },
    @{ Cell = "D46"; Value = "0.7097" },
    @{ Cell = "E46"; Value = "-63.21%" },
    @{ Cell = "E47"; Value = "-0.09%" }
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    $r.NumberFormat = "@"
    $r.Value = $u.Value
    $r.NumberFormat = "General"
}
